# Generate Report for Handoff
# Update the source-file GUID/hash identifiers and handoff timestamps that
# are regenerated each time the localization status report is produced.

$wb = $excel.ActiveWorkbook

$oldGuid = "d5f9d002-9b30-4d5b-b08a-d45fb2baa9e8"
$newGuid = "cc3c47e8-f237-44aa-8aa4-0a2666c9f915"
$oldHash = "23f8d3c55a7ba1c89ab763a8bdbf5631e7b01b1d"
$newHash = "c6823eaa50926bbf5f91706ecb069bc0b14adc3e"

$newMd = "$newGuid.md"
$newZhXlf = "$newGuid.$newHash.zh-cn.xlf"
$newDeXlf = "$newGuid.$newHash.de-de.xlf"

$newZhTime = "2016-03-09 03:36:38"
$newDeTime = "2016-03-09 03:36:49"

# --- Overview sheet: A2 holds the source markdown file name/hyperlink ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMd

# --- zh-cn sheet: A2 markdown file, C2 handoff xlf file, D2 handoff datetime ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newMd
$wsZh.Range("C2").Value = $newZhXlf
$wsZh.Range("D2").Value = $newZhTime

# --- de-de sheet: A2 markdown file, C2 handoff xlf file, D2 handoff datetime ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newMd
$wsDe.Range("C2").Value = $newDeXlf
$wsDe.Range("D2").Value = $newDeTime
